# Weekly fruit/vegetable price update: a new weekly record is inserted
# as row 310 (pushing the existing rows 310-323 down to 311-324), and the
# displaced data lands in the new row 324.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 310 - this shifts rows
# 310:323 down to 311:324 (and the used range grows to A1:R324).
$ws.Rows("310:310").Insert()

# Populate the newly inserted row with this week's record.
$ws.Cells.Item(310, 1).Value = 10
$ws.Cells.Item(310, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(310, 3).Value = "La Araucanía"
$ws.Cells.Item(310, 4).Value = 45267
$ws.Cells.Item(310, 5).Value = 9
$ws.Cells.Item(310, 6).Value = 100112012
$ws.Cells.Item(310, 7).Value = "Espinaca"
$ws.Cells.Item(310, 8).Value = "Sin especificar"
$ws.Cells.Item(310, 9).Value = "Primera"
$ws.Cells.Item(310, 10).Value = 40
$ws.Cells.Item(310, 11).Value = 10000
$ws.Cells.Item(310, 12).Value = 10000
$ws.Cells.Item(310, 13).Value = 10000
$ws.Cells.Item(310, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(310, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(310, 16).Value = 833
$ws.Cells.Item(310, 17).Value = 12
$ws.Cells.Item(310, 18).Value = "Hortaliza"
